# DemoData.xlsx - "digital desk code added"
# Rename Category -> Login, populate login credentials, add a new
# "Manual Claim" sheet with the claim-intake field layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "Category" -> "Login" and rebuild its content
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item("Category")
$login.Name = "Login"

# Clear out the old "Category Name / Primus_test" content first.
$login.Range("A1:B10").ClearContents()

$login.Range("A1").Value = "UserName"
$login.Range("B1").Value = "Password"

$login.Range("A2").Value = "vishal.singh@primussoft.com"
$login.Range("B2").Value = "Primus1234$"

$login.Range("A3").Value = "vishal.singh@primussoft.com"
$login.Range("B3").Value = "Primus1234"

$login.Hyperlinks.Add($login.Range("A2"), "mailto:vishal.singh@primussoft.com")
$login.Range("A2").Style = "Hyperlink"
$login.Hyperlinks.Add($login.Range("A3"), "mailto:vishal.singh@primussoft.com")
$login.Range("A3").Style = "Hyperlink"

$login.Columns.Item(1).ColumnWidth = 31.21875

# ---------------------------------------------------------------------
# 2) Insert the new "Manual Claim" sheet right after "Login"
# ---------------------------------------------------------------------
$claim = $wb.Worksheets.Add($null, $login)
$claim.Name = "Manual Claim"

# Shared-string insertion order follows how the sheet was actually
# authored: headers A1:S1 first, then some row-2 data, then the
# T1:W1 "onsite contact" headers, their row-2 data, and finally J2.
$headersFirst = @(
  "A", "Loss Date ";
  "B", "Recieved Date";
  "C", "Carrier Name";
  "D", "Carrier Ref No#";
  "E", "Claim Number";
  "F", "Cat Code";
  "G", "Peril Type";
  "H", "Property Type";
  "I", "Claim Source";
  "J", "Policy Number";
  "K", "Policy Holder's Name";
  "L", "Policy Holder's Mobile";
  "M", "Policy Holder's Email";
  "N", "Address Line 1";
  "O", "Address Line 2";
  "P", "City";
  "Q", "ZIP";
  "R", "Country";
  "S", "State"
)

for ($i = 0; $i -lt $headersFirst.Count; $i += 2) {
  $col = $headersFirst[$i]
  $text = $headersFirst[$i + 1]
  $cell = $claim.Range($col + "1")
  $cell.Value = $text
  $cell.Font.Bold = $true
}

$claim.Range("A2").Value = 2021
$claim.Range("B2").Value = 2021
$claim.Range("D2").Value = 808080
$claim.Range("E2").Value = 808084
$claim.Range("K2").Value = "John"
$claim.Range("L2").Value = 9867543212
$claim.Range("M2").Value = "kapil12@mailinator.com"
$claim.Range("N2").Value = "4th Avenue"
$claim.Range("P2").Value = "new york"
$claim.Range("Q2").Value = 12345

$headersSecond = @(
  "T", "Onsite Contact Name";
  "U", "Onsite Contact Mobile";
  "V", "Onsite Contact Email"
)
for ($i = 0; $i -lt $headersSecond.Count; $i += 2) {
  $col = $headersSecond[$i]
  $text = $headersSecond[$i + 1]
  $cell = $claim.Range($col + "1")
  $cell.Value = $text
  $cell.Font.Bold = $true
  $cell.WrapText = $true
}

$claim.Range("T2").Value = "kapil"
$claim.Range("U2").Value = 9876543212
$claim.Range("V2").Value = "kapil@gmail.com"

$claim.Range("W1").Value = "Notes"
$claim.Range("W1").Font.Bold = $true
$claim.Range("W1").WrapText = $true

$claim.Range("W2").Value = "This is test remark"

$claim.Range("J2").Value = "DD989438"

$claim.Rows.Item(1).RowHeight = 28.8

$claim.Range("A2,B2,D2,E2,L2,Q2,U2").HorizontalAlignment = -4131

$claim.Hyperlinks.Add($claim.Range("M2"), "mailto:kapil12@mailinator.com")
$claim.Range("M2").Style = "Hyperlink"
$claim.Hyperlinks.Add($claim.Range("V2"), "mailto:kapil@gmail.com")
$claim.Range("V2").Style = "Hyperlink"

$claim.Range("A3").Value = 26
$claim.Range("B3").Value = 1
$claim.Range("A4").Value = 1
$claim.Range("B4").Value = 29
$claim.Range("A3:B4").HorizontalAlignment = -4131

$widths = @(
  1, 15.5546875;
  2, 17.5546875;
  3, 14.5546875;
  4, 15.6640625;
  5, 13.33203125;
  6, 14.44140625;
  7, 10.6640625;
  8, 13.21875;
  9, 12.33203125;
  10, 15.88671875;
  11, 18.6640625;
  12, 19.5546875;
  13, 20.88671875;
  14, 17.44140625;
  15, 13.6640625;
  18, 10.88671875;
  20, 14.44140625;
  21, 17.44140625;
  22, 18.33203125;
  23, 16.44140625
)
for ($i = 0; $i -lt $widths.Count; $i += 2) {
  $claim.Columns.Item($widths[$i]).ColumnWidth = $widths[$i + 1]
}

$claim.PageSetup.Orientation = 1

$claim.Range("E2").Select()

$claim.Activate()

$wb.Windows.Item(1).Left = -108
$wb.Windows.Item(1).Top = -108
